# Apply updated cryptos list values (generated from diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "35.287.14"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.885.62"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.64%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "246.23"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.78%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.689"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  -0.71%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "43.33"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +5.50%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.354"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.92%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "53.71"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  -1.90%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0971"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  +1.91%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.161.43"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.71%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.755"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "4.88"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.883.73"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.06%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "35.431.72"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "72.94"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.13%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0₃0820"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.81%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "244.22"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  -1.79%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.95"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.28%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +10.44%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  -6.29%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "165.94"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.49"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.36%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "18.30"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("E30").Value = "  -2.22%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.128.46"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +10.59%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.28"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("E34").Value = "  -4.39%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.15"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("E37").Value = "  -11.74%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.847"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -2.53%  "
$ws.Range("E40").Value = "  +7.21%  "
$ws.Range("E41").Value = "  +2.86%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "17.20"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "96.46"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -5.47%  "
$ws.Range("E44").Value = "  -2.36%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.298.71"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("E47").Value = "  +7.51%  "
$ws.Range("B48").Value = "Gas"
$ws.Range("C48").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "12.31"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.38%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  -5.63%  "
